$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws | Get-Member | Out-String | Write-Host
